$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.455.75"
$ws.Range("E2").Value = "  -1.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.324.63"
$ws.Range("E3").Value = "  -2.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.71"
$ws.Range("E5").Value = "  -2.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.46"
$ws.Range("E6").Value = "  -7.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  -2.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.318.87"
$ws.Range("E9").Value = "  -2.10%  "

$ws.Range("E10").Value = "  -4.74%  "

$ws.Range("E11").Value = "  -2.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.29"
$ws.Range("E12").Value = "  -4.87%  "

$ws.Range("E13").Value = "  -4.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "667.10"
$ws.Range("E14").Value = "  +4.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.864.08"
$ws.Range("E15").Value = "  -1.84%  "

$ws.Range("E16").Value = "  -2.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.594.77"
$ws.Range("E17").Value = "  -1.65%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.337.97"
$ws.Range("E18").Value = "  -1.96%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.118"
$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.39"
$ws.Range("E20").Value = "  -3.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.90"
$ws.Range("E21").Value = "  -2.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.887"
$ws.Range("E22").Value = "  -2.98%  "

$ws.Range("E23").Value = "  +5.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.95"
$ws.Range("E24").Value = "  -6.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.46"
$ws.Range("E25").Value = "  -1.70%  "

$ws.Range("E26").Value = "  -6.68%  "

$ws.Range("E27").Value = "  -6.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.76"
$ws.Range("E28").Value = "  +2.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.21"
$ws.Range("E29").Value = "  -6.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.39"
$ws.Range("E30").Value = "  -3.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.34"
$ws.Range("E31").Value = "  +6.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "592.25"
$ws.Range("E32").Value = "  -3.72%  "

$ws.Range("E33").Value = "  -2.35%  "

$ws.Range("E34").Value = "  -2.25%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.688.79"
$ws.Range("E36").Value = "  -8.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.73"
$ws.Range("E37").Value = "  -0.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("E38").Value = "  -14.96%  "

$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.08"
$ws.Range("E40").Value = "  -2.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.62"
$ws.Range("E41").Value = "  -6.92%  "

$ws.Range("E42").Value = "  -6.71%  "

$ws.Range("E43").Value = "  -3.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0659"
$ws.Range("E44").Value = "  -7.07%  "

$ws.Range("E45").Value = "  -5.22%  "

$ws.Range("E46").Value = "  -4.83%  "

$ws.Range("E47").Value = "  -1.27%  "

$ws.Range("E48").Value = "  -2.04%  "

$ws.Range("E49").Value = "  -0.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.34"
$ws.Range("E50").Value = "  -3.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "126.72"
$ws.Range("E51").Value = "  -2.52%  "
